# Fill in the homework scores (ДЗ_1..ДЗ_3) for student row 31 (Ялунин Максим)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = 5

# Reflect the cursor position the author ended up at after entering the data
$ws.Range("F31").Select()

$wb.Save()
